# Weekly fruit/vegetable price update: two new "Alcachofa" (artichoke)
# price records for Vega Modelo de Temuco were added to the data table,
# inserted right after the existing row for the prior period (old row 277),
# pushing the rest of the table (old rows 278-317) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 278 (shifts old rows 278-317 down to 280-319).
$ws.Rows.Item(278).Insert()
$ws.Rows.Item(278).Insert()

# --- New row 278 ---
$ws.Cells.Item(278, 1).Value = 10
$ws.Cells.Item(278, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(278, 3).Value = "La Araucanía"
$ws.Cells.Item(278, 4).Value = 45127
$ws.Cells.Item(278, 5).Value = 9
$ws.Cells.Item(278, 6).Value = 100112013
$ws.Cells.Item(278, 7).Value = "Alcachofa"
$ws.Cells.Item(278, 8).Value = "Española"
$ws.Cells.Item(278, 9).Value = "Primera"
$ws.Cells.Item(278, 10).Value = 250
$ws.Cells.Item(278, 11).Value = 17000
$ws.Cells.Item(278, 12).Value = 17000
$ws.Cells.Item(278, 13).Value = 17000
$ws.Cells.Item(278, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(278, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(278, 16).Value = 567
$ws.Cells.Item(278, 17).Value = 30
$ws.Cells.Item(278, 18).Value = "Hortaliza"

# --- New row 279 ---
$ws.Cells.Item(279, 1).Value = 10
$ws.Cells.Item(279, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(279, 3).Value = "La Araucanía"
$ws.Cells.Item(279, 4).Value = 45127
$ws.Cells.Item(279, 5).Value = 9
$ws.Cells.Item(279, 6).Value = 100112013
$ws.Cells.Item(279, 7).Value = "Alcachofa"
$ws.Cells.Item(279, 8).Value = "Madrigal"
$ws.Cells.Item(279, 9).Value = "Primera"
$ws.Cells.Item(279, 10).Value = 300
$ws.Cells.Item(279, 11).Value = 12000
$ws.Cells.Item(279, 12).Value = 12000
$ws.Cells.Item(279, 13).Value = 12000
$ws.Cells.Item(279, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(279, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(279, 16).Value = 300
$ws.Cells.Item(279, 17).Value = 40
$ws.Cells.Item(279, 18).Value = "Hortaliza"
